$p = $ppt.ActivePresentation
$d = $p.Designs.Item(1)
for ($i=1; $i -le $p.Slides.Count; $i++) {
    $p.Slides.Item($i).Design = $d
}
